# Update the trips.xlsx export template:
#  - Replace the Joda-Time based date formatting formulas with the new
#    dateTool.format(...) based formulas.
#  - Move the active selection from D9 to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B6: combined "from - to" period string (row below "Period:" label)
$ws.Cells.Item(6, 2).Value2 = '${dateTool.format("YYYY-MM-dd HH:mm:ss", from, locale, timezone)+" - "+dateTool.format("YYYY-MM-dd HH:mm:ss", to, locale, timezone)}'

# A9: per-trip start time formatted string
$ws.Cells.Item(9, 1).Value2 = '${dateTool.format("YYYY-MM-dd HH:mm:ss", trip.startTime, locale, timezone)}'

# C9: per-trip end time formatted string
$ws.Cells.Item(9, 3).Value2 = '${dateTool.format("YYYY-MM-dd HH:mm:ss", trip.endTime, locale, timezone)}'

# Move the sheet's active selection to B2.
$ws.Range("B2").Select()
